# Update gh-pages output — 南宁-漫展信息.xlsx
# - bump a few "want to go" counters
# - add the newly scraped Canon concert (演出/Performance) event to the
#   "演出" sheet and to the combined "全部类型" sheet (where it is inserted
#   ahead of the existing "第一届ANE·DACG动漫嘉年华" row, pushing that row down)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# 1. 展览 (Exhibition) — simple counter bumps, no structural changes
# ---------------------------------------------------------------------
$wsExpo.Cells.Item(2, 6).Value = 9049   # F2 9015 -> 9049
$wsExpo.Cells.Item(4, 6).Value = 457    # F4 452 -> 457
$wsExpo.Cells.Item(5, 6).Value = 452    # F5 451 -> 452

# ---------------------------------------------------------------------
# 2. 演出 (Performance) — brand-new sheet content, add row 2
# ---------------------------------------------------------------------
# Reuse the bold/bordered/centered style already used for column A of the
# numbering column on the other sheets (cell style index shared with A2 on
# 展览) instead of re-deriving a brand-new style.
$wsExpo.Range("A2").Copy() | Out-Null
$wsShow.Range("A2").PasteSpecial(-4122) | Out-Null

$wsShow.Cells.Item(2, 1).Value = 1
$wsShow.Range("B2:E2").NumberFormat = "@"   # keep these as literal text, not auto-converted dates
$wsShow.Cells.Item(2, 2).Value = "2024-03-30"
$wsShow.Cells.Item(2, 3).Value = "南宁·卡农·世界经典音乐之旅音乐会"
$wsShow.Cells.Item(2, 4).Value = "龙堤路25号 南宁文化艺术中心"
$wsShow.Cells.Item(2, 5).Value = "2024.03.30 20:00-03.30 21:30"
$wsShow.Cells.Item(2, 6).Value = 0
$wsShow.Cells.Item(2, 7).Value = 60
$wsShow.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81798"
$wsShow.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/Tv5lqcVn1707214065277.jpeg"

# ---------------------------------------------------------------------
# 3. 本地生活 (Local life) — untouched
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4. 全部类型 (All types) — same counter bumps as 展览, plus the new
#    Canon concert row inserted at row 5 (the old ANE·DACG row moves to
#    row 6, its counter also bumped to match 展览's new value).
# ---------------------------------------------------------------------
$wsAll.Cells.Item(2, 6).Value = 9049   # F2 9015 -> 9049
$wsAll.Cells.Item(4, 6).Value = 457    # F4 452 -> 457

# Save the current (pre-edit) row 5 contents (ANE·DACG) before overwriting,
# so they can be written back out to the new row 6.
$oldA5 = $wsAll.Cells.Item(5, 1).Value2
$oldB5 = $wsAll.Cells.Item(5, 2).Value2
$oldC5 = $wsAll.Cells.Item(5, 3).Value2
$oldD5 = $wsAll.Cells.Item(5, 4).Value2
$oldE5 = $wsAll.Cells.Item(5, 5).Value2
$oldG5 = $wsAll.Cells.Item(5, 7).Value2
$oldH5 = $wsAll.Cells.Item(5, 8).Value2
$oldI5 = $wsAll.Cells.Item(5, 9).Value2

# Build row 6 first (style for column A copied from the existing numbering
# column so it matches the rest of the sheet).
$wsAll.Range("A5").Copy() | Out-Null
$wsAll.Range("A6").PasteSpecial(-4122) | Out-Null

$wsAll.Cells.Item(6, 1).Value = 5
$wsAll.Range("B6:E6").NumberFormat = "@"
$wsAll.Cells.Item(6, 2).Value = $oldB5
$wsAll.Cells.Item(6, 3).Value = $oldC5
$wsAll.Cells.Item(6, 4).Value = $oldD5
$wsAll.Cells.Item(6, 5).Value = $oldE5
$wsAll.Cells.Item(6, 6).Value = 452     # F6: old ANE-DACG counter, bumped 451 -> 452
$wsAll.Cells.Item(6, 7).Value = $oldG5
$wsAll.Cells.Item(6, 8).Value = $oldH5
$wsAll.Cells.Item(6, 9).Value = $oldI5

# Now overwrite row 5 with the new Canon concert details (A5/G5 stay as-is).
$wsAll.Range("C5:E5").NumberFormat = "@"
$wsAll.Cells.Item(5, 3).Value = "南宁·卡农·世界经典音乐之旅音乐会"
$wsAll.Cells.Item(5, 4).Value = "龙堤路25号 南宁文化艺术中心"
$wsAll.Cells.Item(5, 5).Value = "2024.03.30 20:00-03.30 21:30"
$wsAll.Cells.Item(5, 6).Value = 0       # F5 451 -> 0
$wsAll.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81798"
$wsAll.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/Tv5lqcVn1707214065277.jpeg"
